$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 7166.6665
$ws.Cells.Item(43, 10).Value = 7166.6665
$ws.Cells.Item(43, 12).Value = 7166.6665
$ws.Cells.Item(43, 14).Value = -7304.6665
$ws.Cells.Item(88, 8).Value = 1807.5714
$ws.Cells.Item(88, 9).Value = 1360.6
$ws.Cells.Item(88, 10).Value = 2055.889
$ws.Cells.Item(88, 11).Value = 1360.6
$ws.Cells.Item(88, 12).Value = 2055.889
$ws.Cells.Item(88, 13).Value = -954.5999999999999
$ws.Cells.Item(88, 14).Value = -2867.889
$ws.Cells.Item(91, 8).Value = 1807.5714
$ws.Cells.Item(91, 9).Value = 1360.6
$ws.Cells.Item(91, 10).Value = 2055.889
$ws.Cells.Item(91, 11).Value = 1360.6
$ws.Cells.Item(91, 12).Value = 2055.889
$ws.Cells.Item(91, 13).Value = 43.40000000000009
$ws.Cells.Item(91, 14).Value = -4863.889
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).ClearContents()
$ws.Cells.Item(121, 14).Value = 0
$ws.Cells.Item(137, 8).Value = 2465.3333
$ws.Cells.Item(137, 9).Value = 2465.3333
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 7395.999899999999
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -4845.999899999999
$ws.Cells.Item(141, 8).Value = 5666.3335
$ws.Cells.Item(141, 10).Value = 5999.5
$ws.Cells.Item(141, 12).Value = 17998.5
$ws.Cells.Item(141, 14).Value = -28358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2834.2144
$ws.Cells.Item(45, 9).Value = 1961.875
$ws.Cells.Item(45, 11).Value = 1961.875
$ws.Cells.Item(45, 13).Value = -1584.875
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).ClearContents()
$ws.Cells.Item(59, 14).Value = 0
$ws.Cells.Item(61, 8).Value = 6114.1113
$ws.Cells.Item(61, 9).Value = 6432.5713
$ws.Cells.Item(61, 11).Value = 6432.5713
$ws.Cells.Item(61, 13).Value = -6220.5713
$ws.Cells.Item(132, 8).Value = 3722.5293
$ws.Cells.Item(132, 9).Value = 3347
$ws.Cells.Item(132, 10).Value = 4943
$ws.Cells.Item(132, 11).Value = 10041
$ws.Cells.Item(132, 12).Value = 14829
$ws.Cells.Item(132, 13).Value = -7511
$ws.Cells.Item(132, 14).Value = -19889
$ws.Cells.Item(136, 8).Value = 6114.1113
$ws.Cells.Item(136, 9).Value = 6432.5713
$ws.Cells.Item(136, 11).Value = 19297.7139
$ws.Cells.Item(136, 13).Value = -16747.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 199.5
$ws.Cells.Item(22, 9).Value = 199.5
$ws.Cells.Item(22, 11).Value = 199.5
$ws.Cells.Item(22, 13).Value = -26.5
$ws.Cells.Item(99, 8).Value = 751
$ws.Cells.Item(99, 9).Value = 751
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 751
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).Value = 747
$ws.Cells.Item(105, 8).Value = 4489.0586
$ws.Cells.Item(105, 9).Value = 4181.7
$ws.Cells.Item(105, 11).Value = 4181.7
$ws.Cells.Item(105, 13).Value = -2434.7
$ws.Cells.Item(107, 8).Value = 946.26666
$ws.Cells.Item(107, 9).Value = 952.38464
$ws.Cells.Item(107, 11).Value = 952.38464
$ws.Cells.Item(107, 13).Value = 967.61536

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 1924775.6
$ws.Cells.Item(6, 9).Value = 1364371.6
$ws.Cells.Item(6, 11).Value = 1364371.6
$ws.Cells.Item(6, 13).Value = -1364258.6
$ws.Cells.Item(16, 8).Value = 2126.4
$ws.Cells.Item(16, 9).Value = 2358
$ws.Cells.Item(16, 10).Value = 1200
$ws.Cells.Item(16, 11).Value = 2358
$ws.Cells.Item(16, 12).Value = 1200
$ws.Cells.Item(16, 13).Value = -2071
$ws.Cells.Item(16, 14).Value = -1774
$ws.Cells.Item(58, 8).Value = 3607.6365
$ws.Cells.Item(58, 9).Value = 3077.111
$ws.Cells.Item(58, 11).Value = 3077.111
$ws.Cells.Item(58, 13).Value = -2874.111
$ws.Cells.Item(99, 8).Value = 12370.934
$ws.Cells.Item(99, 9).Value = 6572.5625
$ws.Cells.Item(99, 11).Value = 6572.5625
$ws.Cells.Item(99, 13).Value = -5074.5625
$ws.Cells.Item(107, 8).Value = 1230.6296
$ws.Cells.Item(107, 9).Value = 957.7
$ws.Cells.Item(107, 10).Value = 1391.1765
$ws.Cells.Item(107, 11).Value = 957.7
$ws.Cells.Item(107, 12).Value = 1391.1765
$ws.Cells.Item(107, 13).Value = 962.3
$ws.Cells.Item(107, 14).Value = -5231.1765
$ws.Cells.Item(113, 8).Value = 2126.4
$ws.Cells.Item(113, 9).Value = 2358
$ws.Cells.Item(113, 10).Value = 1200
$ws.Cells.Item(113, 11).Value = 2358
$ws.Cells.Item(113, 12).Value = 1200
$ws.Cells.Item(113, 13).Value = -188
$ws.Cells.Item(113, 14).Value = -5540
$ws.Cells.Item(122, 8).Value = 2994
$ws.Cells.Item(122, 9).Value = 2994
$ws.Cells.Item(122, 11).Value = 8982
$ws.Cells.Item(122, 13).Value = -6532
$ws.Cells.Item(126, 8).Value = 12370.934
$ws.Cells.Item(126, 9).Value = 6572.5625
$ws.Cells.Item(126, 11).Value = 19717.6875
$ws.Cells.Item(126, 13).Value = -17247.6875
$ws.Cells.Item(132, 8).Value = 9220
$ws.Cells.Item(132, 9).Value = 9997.5
$ws.Cells.Item(132, 10).Value = 8598
$ws.Cells.Item(132, 11).Value = 29992.5
$ws.Cells.Item(132, 12).Value = 25794
$ws.Cells.Item(132, 13).Value = -27462.5
$ws.Cells.Item(132, 14).Value = -30854
$ws.Cells.Item(136, 8).Value = 3607.6365
$ws.Cells.Item(136, 9).Value = 3077.111
$ws.Cells.Item(136, 11).Value = 9231.332999999999
$ws.Cells.Item(136, 13).Value = -6681.332999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 127.15
$ws.Cells.Item(12, 9).Value = 130.22223
$ws.Cells.Item(12, 11).Value = 390.66669
$ws.Cells.Item(12, 13).Value = -217.66669
$ws.Cells.Item(88, 8).Value = 3250
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(91, 8).Value = 3250
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 2417.8
$ws.Cells.Item(132, 9).Value = 2417.8
$ws.Cells.Item(132, 11).Value = 21760.2
$ws.Cells.Item(132, 13).Value = -19230.2
$ws.Cells.Item(137, 8).Value = 3529.2222
$ws.Cells.Item(137, 9).Value = 3594.6
$ws.Cells.Item(137, 11).Value = 10783.8
$ws.Cells.Item(137, 13).Value = -5683.799999999999
$ws.Cells.Item(141, 8).Value = 8255.75
$ws.Cells.Item(141, 9).Value = 8255.75
$ws.Cells.Item(141, 11).Value = 24767.25
$ws.Cells.Item(141, 13).Value = -19587.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 9
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 9
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).ClearContents()
$ws.Cells.Item(13, 13).Value = 9
$ws.Cells.Item(13, 14).Value = -287
$ws.Cells.Item(36, 8).Value = 9499.5
$ws.Cells.Item(36, 9).Value = 11799.4
$ws.Cells.Item(36, 11).Value = 11799.4
$ws.Cells.Item(36, 13).Value = -11314.4
$ws.Cells.Item(53, 8).Value = 30000.5
$ws.Cells.Item(53, 10).Value = 30000.5
$ws.Cells.Item(53, 12).Value = 30000.5
$ws.Cells.Item(53, 14).Value = -31262.5
$ws.Cells.Item(80, 8).Value = 3998.8
$ws.Cells.Item(80, 9).Value = 3998.5
$ws.Cells.Item(80, 11).Value = 3998.5
$ws.Cells.Item(80, 13).Value = -3000.5
$ws.Cells.Item(83, 8).Value = 3998.8
$ws.Cells.Item(83, 9).Value = 3998.5
$ws.Cells.Item(83, 11).Value = 19992.5
$ws.Cells.Item(83, 13).Value = -15000.5
$ws.Cells.Item(102, 8).Value = 1112.2142
$ws.Cells.Item(102, 9).Value = 961.0909
$ws.Cells.Item(102, 11).Value = 961.0909
$ws.Cells.Item(102, 13).Value = 660.9091
$ws.Cells.Item(122, 8).Value = 204101.2
$ws.Cells.Item(122, 9).Value = 4502.3335
$ws.Cells.Item(122, 10).Value = 503499.5
$ws.Cells.Item(122, 11).Value = 13507.0005
$ws.Cells.Item(122, 12).Value = 1510498.5
$ws.Cells.Item(122, 13).Value = -11057.0005
$ws.Cells.Item(122, 14).Value = -1515398.5
$ws.Cells.Item(126, 8).Value = 5829.875
$ws.Cells.Item(126, 9).Value = 5551.25
$ws.Cells.Item(126, 11).Value = 16653.75
$ws.Cells.Item(126, 13).Value = -14183.75
$ws.Cells.Item(132, 8).Value = 2285.889
$ws.Cells.Item(132, 9).Value = 2305.4856
$ws.Cells.Item(132, 10).Value = 1600
$ws.Cells.Item(132, 11).Value = 6916.4568
$ws.Cells.Item(132, 12).Value = 4800
$ws.Cells.Item(132, 13).Value = -4386.4568
$ws.Cells.Item(132, 14).Value = -9860

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 14899
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 14899
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).ClearContents()
$ws.Cells.Item(16, 13).Value = 14899
$ws.Cells.Item(16, 14).Value = -15239
$ws.Cells.Item(20, 8).Value = 63797.188
$ws.Cells.Item(20, 9).Value = 92436.82000000001
$ws.Cells.Item(20, 10).Value = 790
$ws.Cells.Item(20, 11).Value = 92436.82000000001
$ws.Cells.Item(20, 12).Value = 790
$ws.Cells.Item(20, 13).Value = -92210.82000000001
$ws.Cells.Item(20, 14).Value = -1242
$ws.Cells.Item(46, 8).Value = 3916.0833
$ws.Cells.Item(46, 9).Value = 2375
$ws.Cells.Item(46, 10).Value = 6998.25
$ws.Cells.Item(46, 11).Value = 2375
$ws.Cells.Item(46, 12).Value = 6998.25
$ws.Cells.Item(46, 13).Value = -2187
$ws.Cells.Item(46, 14).Value = -7374.25
$ws.Cells.Item(100, 8).Value = 2269.2
$ws.Cells.Item(100, 9).Value = 2274.125
$ws.Cells.Item(100, 11).Value = 2274.125
$ws.Cells.Item(100, 13).Value = -1733.125
$ws.Cells.Item(122, 8).Value = 6243.75
$ws.Cells.Item(122, 9).Value = 4359
$ws.Cells.Item(122, 11).Value = 13077
$ws.Cells.Item(122, 13).Value = -10627

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 75000
$ws.Cells.Item(46, 10).Value = 75000
$ws.Cells.Item(46, 12).Value = 75000
$ws.Cells.Item(46, 14).Value = -75462
$ws.Cells.Item(48, 8).Value = 45683.332
$ws.Cells.Item(48, 9).Value = 50000
$ws.Cells.Item(48, 11).Value = 50000
$ws.Cells.Item(48, 13).Value = -49431
$ws.Cells.Item(81, 8).Value = 1424.8334
$ws.Cells.Item(81, 9).Value = 1424.8334
$ws.Cells.Item(81, 11).Value = 2849.6668
$ws.Cells.Item(81, 13).Value = -1788.6668
$ws.Cells.Item(84, 8).Value = 1424.8334
$ws.Cells.Item(84, 9).Value = 1424.8334
$ws.Cells.Item(84, 11).Value = 14248.334
$ws.Cells.Item(84, 13).Value = -8944.333999999999
$ws.Cells.Item(122, 8).Value = 2383.04
$ws.Cells.Item(122, 9).Value = 1782.6111
$ws.Cells.Item(122, 10).Value = 3927
$ws.Cells.Item(122, 11).Value = 5347.8333
$ws.Cells.Item(122, 12).Value = 11781
$ws.Cells.Item(122, 13).Value = -2897.8333
$ws.Cells.Item(122, 14).Value = -16681
$ws.Cells.Item(134, 8).Value = 75000
$ws.Cells.Item(134, 10).Value = 75000
$ws.Cells.Item(134, 12).Value = 225000
$ws.Cells.Item(134, 14).Value = -230070
$ws.Cells.Item(136, 8).Value = 3324.8
$ws.Cells.Item(136, 9).Value = 3297.923
$ws.Cells.Item(136, 11).Value = 9893.769
$ws.Cells.Item(136, 13).Value = -7343.769
